# Budget breakdown: split "Scientific Personnel" into two line items
# (Research Assistant + PhD student), add description/amount columns for
# several other cost categories, and total everything with a SUM formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Insert a new row after row 5 so the "Scientific Personnel" category gets
# a second line (row 6), pushing all the following categories down by one.
$ws.Rows.Item(6).Insert()

# Row 5: Scientific Personnel - Research Assistant
# (leading "'" forces a text/quote-prefixed cell, matching the style Excel
# recorded for these two new description cells; it is stripped from the
# stored value itself.)
$ws.Range("B5").Value = "'One Research Assistant to work on methodological development (1 year)"
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = 23000

# Row 6: Scientific Personnel - PhD student (A6 left blank; category label stays on row 5)
$ws.Range("B6").Value = "'One PhD student to conduct empirical research (18 months)"
$ws.Range("B6").WrapText = $true
$ws.Range("C6").Value = 67500

# Row 7: Other Personnel
$ws.Range("C7").Value = 0

# Row 8: Travel Costs
$ws.Range("B8").Value = "Travel to conferences and workshops for PI and research assistant. Potential visit to the Demography Unit of Stockholm University (Sweden) to work on micro-level data"
$ws.Range("C8").Value = 7000

# Row 9: Publication costs
$ws.Range("B9").Value = "Article processing and open access charges (expected a minimum of three publications in top journals)"
$ws.Range("C9").Value = 6000

# Row 10: Other Recurring Costs
$ws.Range("C10").Value = 0

# Row 11: Equipment
$ws.Range("C11").Value = 0

# Row 12: Other Non-recurring Costs
$ws.Range("C12").Value = 0

# Row 13: Total
$ws.Range("C13").Formula = "=SUM(C5:C12)"

# Restore the view state recorded in the saved file.
$ws.Range("I10").Select()
